$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most Price (D) values look like plain numbers (e.g. "560.56"); Excel would
# auto-convert those to numeric cells, but the sheet stores Price as text, so
# for such values we force the cell format to Text before assigning.

$ws.Range("D2").Value = "68.885.81"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "2.441.20"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.56"
$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.41"
$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  +0.82%  "

$ws.Range("E9").Value = "  +11.22%  "

$ws.Range("E10").Value = "  -1.78%  "

$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.56"
$ws.Range("E12").Value = "  -5.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000177"
$ws.Range("E13").Value = "  +5.87%  "

$ws.Range("D14").Value = "68.792.79"
$ws.Range("E14").Value = "  +0.31%  "

$ws.Range("D15").Value = "2.887.96"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.31"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "2.442.66"
$ws.Range("E17").Value = "  -0.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.56"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.16"
$ws.Range("E19").Value = "  +1.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.95"
$ws.Range("E20").Value = "  +0.87%  "

$ws.Range("E21").Value = "  +2.61%  "

$ws.Range("E22").Value = "  +2.86%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.12"
$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("E25").Value = "  +2.88%  "

$ws.Range("D26").Value = "2.569.82"
$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.22"
$ws.Range("E28").Value = "  +1.17%  "

$ws.Range("D29").Value = "0.0₃0824"
$ws.Range("E29").Value = "  +1.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.10%  "

# Coinranking shuffled positions 31/32: Bittensor and Fetch.AI swapped rows
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("E32").Value = "  +3.12%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "429.83"
$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.60"
$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.92"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.99"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.00"
$ws.Range("E38").Value = "  +1.50%  "

$ws.Range("E39").Value = "  -1.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.299"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.51"
$ws.Range("E41").Value = "  +4.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.36"
$ws.Range("E42").Value = "  -0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.08"
$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.05"
$ws.Range("E44").Value = "  +0.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.35"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.10"
$ws.Range("E46").Value = "  +0.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0718"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.481"
$ws.Range("E48").Value = "  +0.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.557"
$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0924"
$ws.Range("E50").Value = "  +1.61%  "

$ws.Range("E51").Value = "  +1.10%  "
